$d = $word.ActiveDocument
$targetPara = $d.Paragraphs.Item(85)
$r = $targetPara.Range

$xmlFragment = @'
<w:p w14:paraId="2F0BC143" w14:textId="77777777" w:rsidR="00A85D2D" w:rsidRDefault="00A85D2D" w:rsidP="00A85D2D">
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t>Video presentación de resolución del trabajo práctico número uno de Análisis Numérico correspondiente al ciclo lectivo 2025.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Para la resolución de este primer trabajo práctico de Análisis Numérico año </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>2025,hemos</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:t xml:space="preserve"> recibido tres archivos que contenían lecturas de señales de electroencefalogramas, corres</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">pondientes a distintas fases del ciclo de la epilepsia, a saber: sano, </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>interictal</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> y </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>convusión</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">La frecuencia de muestreo de las señales fue informada: 173.61 </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>herzios</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Para analizar las señales, aplicamos un filtro pasa bajos en 40 </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>herzios</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>, frecuencia límite en la cual la señal del electroencefalograma no está contaminada por ruidos musculares, ni aparatos eléctricos del entorno.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Para llevar las lecturas desde el dominio del tiempo al dominio de las frecuencias, utilizamos la transformada rápida de </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>fourier</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Hicimos </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>una análisis</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:t xml:space="preserve"> luego de la potencia espectral, lo que nos permitió determinar en qué </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>frencuencias</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> se concentraba más energía en la señal, y hacer una comparativa entre la fase de “sano” versus la </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>interictal</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t xml:space="preserve"> y la convulsión.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t xml:space="preserve">Finalmente, calculamos la autocorrelación. Al ser una operación de </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>convolusión</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>, llegando a cero tiende a ser máxima, porque allí es donde todos los valores confluyen.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Las dificultades que se nos presentaron al momento de la resolución estuvieron vinculadas a cuestiones organizativas, por una parte, y técnicas por la otra. Aspectos que fueron resueltos, creemos, </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:t>postivamente</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:t>.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t xml:space="preserve">En cuanto a los aspectos destacados para la resolución, la herramienta fundamental fue la transformada de Fourier en tiempo discreto, utilizando la librería </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t>Scipy</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t xml:space="preserve"> provista por el lenguaje Python.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t xml:space="preserve">Conforme lo planteado al inicio, y gracias al uso del </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t>pasabajos</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t xml:space="preserve">, las conclusiones arribadas permiten afirmar qué: en el ciclo sano, la energía se concentra en frecuencias más bajas; en el </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t>interictal</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t xml:space="preserve"> la distribución es dispersa, y en la convulsión la distribución es caótica.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t>La densidad espectral de potencia mostró que hay mayor energía en frecuencias inusuales respecto a la fase sano.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t>La autocorrelación mostró que la fase sana presenta patrones periódicos.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t>Muchas gracias por su atención.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:ind w:firstLine="720"/>
  </w:pPr>
</w:p>
'@

$r.InsertXML($xmlFragment)

Write-Output "Paragraphs count: $($d.Paragraphs.Count)"
